# v1.2.3.1 P! - Js e CSS das Contas OK
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Fill in the missing Status for row 18, and swap the descriptions of
# rows 18 and 19 (G18 <-> G19).
$ws.Range("B18").Value = "Ok"
$ws.Range("G18").Value = "Js Tela - Controle de Contas"
$ws.Range("G19").Value = "Alterar Transferencia de Contas"

# Move the active selection to B19
$ws.Range("B19").Select()
